$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 updates (day 103 post added; shift the "ser" references)
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 103"
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 101"
$ws.Range("I7").Value = 100
